$d = $word.ActiveDocument

# --- locate the paragraph that ends with
# "... out of the scope of this project." -----------------------------
$needle = "out of the scope of this project."
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    $needle, $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target text: '$needle'"
}

$targetParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Start -le $searchRange.Start -and
        $candidate.Range.End -ge $searchRange.End) {
        $targetParagraph = $candidate
        break
    }
}

if ($null -eq $targetParagraph) {
    throw "Could not locate the paragraph containing the target text."
}

# --- insert a brand-new paragraph right after it ----------------------
# (its paragraph/run formatting - justify "both" + Times New Roman - is
# inherited automatically from the paragraph it was split off from,
# same as the surrounding paragraphs already in the document)
$tailRange = $targetParagraph.Range
$tailRange.Collapse(0)            # wdCollapseEnd
$insertPos = $tailRange.Start
$tailRange.InsertParagraphAfter()

# --- fill the freshly-created (still empty) paragraph with the two
# new sentences ---------------------------------------------------------
$firstSentence = "As of 2026, F1 is planning to use 100% fully sustainable fuels. For this study this fuel is assumed to be ethanol as "
$secondSentence = "series such as BTCC have already started using ethanol as of 2025."

$newRange = $d.Range($insertPos, $insertPos)
$newRange.Expand(4) | Out-Null    # wdParagraph - grabs the whole new paragraph
$newRange.InsertBefore($firstSentence)

# move past the text we just typed (still inside the paragraph, ahead
# of its end-of-paragraph mark) before appending the second sentence
$newRange.MoveStart(1, $firstSentence.Length) | Out-Null
$newRange.InsertBefore($secondSentence)
